# Apply crypto price/volume updates as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.558.56"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "2.595.43"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'516.81"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").Value = "'153.38"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.595"
$ws.Range("E8").Value = "  +2.59%  "
$ws.Range("D9").Value = "'6.67"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D13").Value = "3.050.89"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").Value = "60.524.00"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "'21.70"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "'0.0000140"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("D17").Value = "2.599.22"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D19").Value = "'351.77"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").Value = "'10.58"
$ws.Range("E20").Value = "  +2.65%  "
$ws.Range("D21").Value = "'6.22"
$ws.Range("E21").Value = "  +2.67%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "'61.05"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("D24").Value = "'0.428"
$ws.Range("E24").Value = "  +2.59%  "
$ws.Range("D25").Value = "2.712.94"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "0.0₃0842"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "'7.33"
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "'6.27"
$ws.Range("E31").Value = "  +9.82%  "
$ws.Range("D32").Value = "'19.41"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("D34").Value = "'150.58"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("E35").Value = "  +3.09%  "
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("D37").Value = "'0.916"
$ws.Range("E37").Value = "  +6.61%  "
$ws.Range("D38").Value = "'1.48"
$ws.Range("E38").Value = "  +2.68%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'36.36"
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'3.78"
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("D41").Value = "'0.837"
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("D42").Value = "'286.83"
$ws.Range("E42").Value = "  -4.28%  "
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("D44").Value = "'0.624"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").Value = "'0.0559"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").Value = "'19.52"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").Value = "'0.0236"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'18.99"
$ws.Range("E51").Value = "  +8.22%  "
